$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44446
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 30000
$ws.Range("S2").Value = 3000
$ws.Range("D3").Value = 44454
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 320
$ws.Range("N3").Value = 30000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 30000
$ws.Range("S3").Value = 3000
$ws.Range("D4").Value = 44454
$ws.Range("M4").Value = 300
$ws.Range("N4").Value = 28000
$ws.Range("O4").Value = 28000
$ws.Range("P4").Value = 28000
$ws.Range("S4").Value = 2800
$ws.Range("D5").Value = 44489
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 27000
$ws.Range("O5").Value = 27000
$ws.Range("P5").Value = 27000
$ws.Range("S5").Value = 2700
$ws.Range("D6").Value = 44489
$ws.Range("M6").Value = 50
$ws.Range("N6").Value = 25000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 25000
$ws.Range("S6").Value = 2500
$ws.Range("D7").Value = 44503
$ws.Range("M7").Value = 140
$ws.Range("D8").Value = 44522
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 25000
$ws.Range("O8").Value = 25000
$ws.Range("P8").Value = 25000
$ws.Range("S8").Value = 2500
$ws.Range("D9").Value = 44522
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 23000
$ws.Range("O9").Value = 23000
$ws.Range("P9").Value = 23000
$ws.Range("S9").Value = 2300
$ws.Range("D10").Value = 44475
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 200
$ws.Range("N10").Value = 28000
$ws.Range("O10").Value = 28000
$ws.Range("P10").Value = 28000
$ws.Range("S10").Value = 2800
$ws.Range("D11").Value = 44510
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 150
$ws.Range("D12").Value = 44511
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 25000
$ws.Range("O12").Value = 25000
$ws.Range("P12").Value = 25000
$ws.Range("S12").Value = 2500
$ws.Range("D14").Value = 44512
$ws.Range("L14").Value = "Especial"
$ws.Range("N14").Value = 26000
$ws.Range("O14").Value = 26000
$ws.Range("P14").Value = 26000
$ws.Range("S14").Value = 2600
$ws.Range("D15").Value = 44462
$ws.Range("L15").Value = "Especial"
$ws.Range("M15").Value = 205
$ws.Range("N15").Value = 30000
$ws.Range("O15").Value = 30000
$ws.Range("P15").Value = 30000
$ws.Range("S15").Value = 3000
$ws.Range("D16").Value = 44462
$ws.Range("M16").Value = 180
$ws.Range("N16").Value = 28000
$ws.Range("O16").Value = 28000
$ws.Range("P16").Value = 28000
$ws.Range("S16").Value = 2800
$ws.Range("D17").Value = 44461
$ws.Range("M17").Value = 150
$ws.Range("D18").Value = 44461
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 25000
$ws.Range("O18").Value = 25000
$ws.Range("P18").Value = 25000
$ws.Range("S18").Value = 2500
$ws.Range("D19").Value = 44467
$ws.Range("L19").Value = "Especial"
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 30000
$ws.Range("O19").Value = 30000
$ws.Range("P19").Value = 30000
$ws.Range("S19").Value = 3000
$ws.Range("D20").Value = 44467
$ws.Range("L20").Value = "Primera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 28000
$ws.Range("O20").Value = 28000
$ws.Range("P20").Value = 28000
$ws.Range("S20").Value = 2800
$ws.Range("D21").Value = 44508
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 50
$ws.Range("N21").Value = 25000
$ws.Range("O21").Value = 25000
$ws.Range("P21").Value = 25000
$ws.Range("S21").Value = 2500
$ws.Range("D22").Value = 44508
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 60
$ws.Range("N22").Value = 23000
$ws.Range("O22").Value = 23000
$ws.Range("P22").Value = 23000
$ws.Range("S22").Value = 2300
$ws.Range("D23").Value = 44459
$ws.Range("L23").Value = "Especial"
$ws.Range("N23").Value = 30000
$ws.Range("O23").Value = 30000
$ws.Range("P23").Value = 30000
$ws.Range("S23").Value = 3000
$ws.Range("D24").Value = 44445
$ws.Range("M24").Value = 250
$ws.Range("N24").Value = 28000
$ws.Range("O24").Value = 30000
$ws.Range("P24").Value = 29200
$ws.Range("S24").Value = 2920
$ws.Range("D25").Value = 44431
$ws.Range("M25").Value = 30
$ws.Range("D26").Value = 44495
$ws.Range("L26").Value = "Primera"
$ws.Range("N26").Value = 25000
$ws.Range("O26").Value = 25000
$ws.Range("P26").Value = 25000
$ws.Range("S26").Value = 2500
$ws.Range("D27").Value = 44523
$ws.Range("M27").Value = 150
$ws.Range("N27").Value = 23000
$ws.Range("O27").Value = 23000
$ws.Range("P27").Value = 23000
$ws.Range("S27").Value = 2300
$ws.Range("D28").Value = 44438
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 100
$ws.Range("D30").Value = 44466
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 110
$ws.Range("N30").Value = 30000
$ws.Range("O30").Value = 30000
$ws.Range("P30").Value = 30000
$ws.Range("S30").Value = 3000
$ws.Range("D31").Value = 44498
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 22000
$ws.Range("O31").Value = 23000
$ws.Range("P31").Value = 22600
$ws.Range("S31").Value = 2260
$ws.Range("D32").Value = 44468
$ws.Range("L32").Value = "Especial"
$ws.Range("M32").Value = 250
$ws.Range("N32").Value = 30000
$ws.Range("O32").Value = 30000
$ws.Range("P32").Value = 30000
$ws.Range("S32").Value = 3000
$ws.Range("D33").Value = 44434
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 30000
$ws.Range("O33").Value = 30000
$ws.Range("P33").Value = 30000
$ws.Range("S33").Value = 3000
$ws.Range("D34").Value = 44530
$ws.Range("L34").Value = "Primera"
$ws.Range("N34").Value = 20000
$ws.Range("O34").Value = 20000
$ws.Range("P34").Value = 20000
$ws.Range("S34").Value = 2000
$ws.Range("D35").Value = 44517
$ws.Range("L35").Value = "Especial"
$ws.Range("M35").Value = 70
$ws.Range("N35").Value = 27000
$ws.Range("O35").Value = 27000
$ws.Range("P35").Value = 27000
$ws.Range("S35").Value = 2700
$ws.Range("D36").Value = 44517
$ws.Range("M36").Value = 80
$ws.Range("N36").Value = 25000
$ws.Range("O36").Value = 25000
$ws.Range("P36").Value = 25000
$ws.Range("S36").Value = 2500
$ws.Range("D39").Value = 44441
$ws.Range("M39").Value = 150
$ws.Range("D40").Value = 44515
$ws.Range("N40").Value = 25000
$ws.Range("O40").Value = 25000
$ws.Range("P40").Value = 25000
$ws.Range("S40").Value = 2500
$ws.Range("D41").Value = 44447
$ws.Range("L41").Value = "Especial"
$ws.Range("M41").Value = 50
$ws.Range("N41").Value = 32000
$ws.Range("O41").Value = 32000
$ws.Range("P41").Value = 32000
$ws.Range("S41").Value = 3200
$ws.Range("D42").Value = 44453
$ws.Range("L42").Value = "Especial"
$ws.Range("M42").Value = 135
$ws.Range("N42").Value = 30000
$ws.Range("O42").Value = 30000
$ws.Range("P42").Value = 30000
$ws.Range("S42").Value = 3000
$ws.Range("D43").Value = 44487
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 130
$ws.Range("N43").Value = 25000
$ws.Range("O43").Value = 25000
$ws.Range("P43").Value = 25000
$ws.Range("S43").Value = 2500
$ws.Range("D44").Value = 44496
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 200
$ws.Range("N44").Value = 25000
$ws.Range("O44").Value = 25000
$ws.Range("P44").Value = 25000
$ws.Range("S44").Value = 2500
$ws.Range("D45").Value = 44490
$ws.Range("M45").Value = 120
$ws.Range("N45").Value = 25000
$ws.Range("O45").Value = 25000
$ws.Range("P45").Value = 25000
$ws.Range("S45").Value = 2500
$ws.Range("D46").Value = 44432
$ws.Range("L46").Value = "Especial"
$ws.Range("M46").Value = 70
$ws.Range("N46").Value = 30000
$ws.Range("O46").Value = 30000
$ws.Range("P46").Value = 30000
$ws.Range("S46").Value = 3000
$ws.Range("D47").Value = 44463
$ws.Range("L47").Value = "Especial"
$ws.Range("M47").Value = 150
$ws.Range("N47").Value = 30000
$ws.Range("O47").Value = 30000
$ws.Range("P47").Value = 30000
$ws.Range("S47").Value = 3000
$ws.Range("D48").Value = 44463
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 100
$ws.Range("N48").Value = 26000
$ws.Range("O48").Value = 26000
$ws.Range("P48").Value = 26000
$ws.Range("S48").Value = 2600
$ws.Range("D49").Value = 44484
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 120
$ws.Range("N49").Value = 25000
$ws.Range("O49").Value = 25000
$ws.Range("P49").Value = 25000
$ws.Range("S49").Value = 2500
$ws.Range("D50").Value = 44484
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = 22000
$ws.Range("O50").Value = 22000
$ws.Range("P50").Value = 22000
$ws.Range("S50").Value = 2200
$ws.Range("D51").Value = 44455
$ws.Range("M51").Value = 150
$ws.Range("N51").Value = 30000
$ws.Range("O51").Value = 30000
$ws.Range("P51").Value = 30000
$ws.Range("S51").Value = 3000
$ws.Range("D52").Value = 44516
$ws.Range("M52").Value = 250
$ws.Range("N52").Value = 25000
$ws.Range("O52").Value = 25000
$ws.Range("P52").Value = 25000
$ws.Range("S52").Value = 2500
$ws.Range("D53").Value = 44491
$ws.Range("M53").Value = 100
$ws.Range("D56").Value = 44482
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 120
$ws.Range("N56").Value = 25000
$ws.Range("O56").Value = 25000
$ws.Range("P56").Value = 25000
$ws.Range("S56").Value = 2500
$ws.Range("D57").Value = 44474
$ws.Range("L57").Value = "Especial"
$ws.Range("M57").Value = 150
$ws.Range("N57").Value = 30000
$ws.Range("O57").Value = 30000
$ws.Range("P57").Value = 30000
$ws.Range("S57").Value = 3000
$ws.Range("D58").Value = 44494
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 150
$ws.Range("N58").Value = 25000
$ws.Range("O58").Value = 25000
$ws.Range("P58").Value = 25000
$ws.Range("S58").Value = 2500
$ws.Range("D59").Value = 44494
$ws.Range("L59").Value = "Segunda"
$ws.Range("M59").Value = 50
$ws.Range("N59").Value = 23000
$ws.Range("O59").Value = 23000
$ws.Range("P59").Value = 23000
$ws.Range("S59").Value = 2300
$ws.Range("D60").Value = 44473
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 200
$ws.Range("N60").Value = 28000
$ws.Range("O60").Value = 28000
$ws.Range("P60").Value = 28000
$ws.Range("S60").Value = 2800
$ws.Range("D61").Value = 44518
$ws.Range("M61").Value = 210
$ws.Range("N61").Value = 20000
$ws.Range("O61").Value = 20000
$ws.Range("P61").Value = 20000
$ws.Range("S61").Value = 2000
$ws.Range("D62").Value = 44505
$ws.Range("M62").Value = 100
$ws.Range("N62").Value = 25000
$ws.Range("O62").Value = 25000
$ws.Range("P62").Value = 25000
$ws.Range("S62").Value = 2500
$ws.Range("D63").Value = 44421
$ws.Range("M63").Value = 30
$ws.Range("N63").Value = 35000
$ws.Range("O63").Value = 35000
$ws.Range("P63").Value = 35000
$ws.Range("S63").Value = 3500
$ws.Range("D64").Value = 44483
$ws.Range("M64").Value = 80
$ws.Range("D65").Value = 44483
$ws.Range("L65").Value = "Segunda"
$ws.Range("M65").Value = 30
$ws.Range("N65").Value = 22000
$ws.Range("O65").Value = 22000
$ws.Range("P65").Value = 22000
$ws.Range("S65").Value = 2200
$ws.Range("D66").Value = 44519
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 120
$ws.Range("N66").Value = 23000
$ws.Range("O66").Value = 23000
$ws.Range("P66").Value = 23000
$ws.Range("S66").Value = 2300
$ws.Range("D67").Value = 44519
$ws.Range("L67").Value = "Segunda"
$ws.Range("M67").Value = 90
$ws.Range("N67").Value = 20000
$ws.Range("O67").Value = 20000
$ws.Range("P67").Value = 20000
$ws.Range("S67").Value = 2000
$ws.Range("D68").Value = 44519
$ws.Range("L68").Value = "Tercera"
$ws.Range("M68").Value = 60
$ws.Range("N68").Value = 17000
$ws.Range("O68").Value = 17000
$ws.Range("P68").Value = 17000
$ws.Range("S68").Value = 1700
